# Update "想去人数" (F column) figures across the four sheets of the
# 杭州-漫展信息 workbook, matching the regenerated-data commit.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12814
$ws1.Range("F3").Value  = 7169
$ws1.Range("F10").Value = 1007
$ws1.Range("F12").Value = 356
$ws1.Range("F16").Value = 1019
$ws1.Range("F18").Value = 249
$ws1.Range("F22").Value = 312
$ws1.Range("F24").Value = 167
$ws1.Range("F25").Value = 373
$ws1.Range("F26").Value = 5239
$ws1.Range("F28").Value = 1431
$ws1.Range("F29").Value = 312
$ws1.Range("F30").Value = 1370
$ws1.Range("F32").Value = 45
$ws1.Range("F36").Value = 596
$ws1.Range("F38").Value = 3736

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 56

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9281
$ws3.Range("F4").Value = 2021

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9281
$ws4.Range("F4").Value  = 2021
$ws4.Range("F5").Value  = 12814
$ws4.Range("F6").Value  = 7169
$ws4.Range("F10").Value = 1007
$ws4.Range("F12").Value = 356
$ws4.Range("F16").Value = 1019
$ws4.Range("F18").Value = 249
$ws4.Range("F22").Value = 312
$ws4.Range("F27").Value = 167
$ws4.Range("F28").Value = 373
$ws4.Range("F29").Value = 5239
$ws4.Range("F31").Value = 1431
$ws4.Range("F34").Value = 312
$ws4.Range("F36").Value = 1370
$ws4.Range("F40").Value = 596
$ws4.Range("F47").Value = 3736
